$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docIds = @{
    2 = "-5024543956076458118"
    3 = "1089035699668342371"
    4 = "-1244263580737547364"
    5 = "-161301453925448720"
    6 = "6718241552717100685"
    7 = "-3905887755394201128"
    8 = "2087618240458196074"
    9 = "-1133448472351489615"
    10 = "-3754667236999779330"
    11 = "1437254590655497329"
    12 = "-3904202706470053090"
    13 = "3765829122322171259"
    14 = "4443689924724281400"
    15 = "342447401836312909"
    16 = "-3748958274366039936"
    17 = "5373755885703083921"
    18 = "234503726462468368"
    19 = "2683659211737244996"
    20 = "2759954696691783152"
    21 = "-5187000563749245230"
    22 = "1011090460398641051"
    23 = "-4103418074422834420"
    24 = "-6376283724736182054"
    25 = "5823057369151609802"
    26 = "-5283959848959794851"
    27 = "-7599315705937794156"
    28 = "-182851223443396917"
    29 = "520941480054210267"
    30 = "8642926451828596316"
    31 = "-7492451031478652905"
    32 = "1477741574240669269"
    33 = "-7240450692410412136"
    34 = "6110845092117184315"
    35 = "-967750167548269064"
    36 = "-6575885249756271861"
    37 = "5331645556016952874"
    38 = "-4262106502305468479"
    39 = "-634420195667353983"
    40 = "-2123348678845688084"
    41 = "3044957725242547638"
    42 = "-1032778900891362717"
    43 = "9078632615318898839"
    44 = "3721174555711123392"
    45 = "-138546740687672117"
    46 = "-4038202933620396616"
    47 = "-7510076755748020832"
    48 = "400801814425465643"
    49 = "-5488338553918754255"
    50 = "-5894097881799099125"
    51 = "2674234454793678402"
    52 = "-2770140877516060079"
    53 = "6286707755240707741"
    54 = "-8908226992892378639"
    55 = "-7477131404435700751"
    56 = "338897088895821016"
    57 = "-1126872237344345472"
    58 = "6862186950629248297"
    59 = "7036815626524267464"
    60 = "-695862899884876166"
    61 = "-2587052032452335265"
    62 = "230596285226370772"
    63 = "-2541033976428744266"
    64 = "-1208576990337673422"
    65 = "4368958120626670240"
    66 = "7720003327581256505"
    67 = "-1887934018914723522"
    68 = "6952037141390866772"
    69 = "2155417342542719627"
    70 = "6539586539588433881"
    71 = "1994887170751869595"
    72 = "4332347659672953139"
    73 = "-8399382642519882028"
    74 = "3965854485100243398"
    75 = "8208960443869166268"
    76 = "-353724325541503250"
    77 = "-4987400653008156332"
    78 = "-5555794890157511271"
    79 = "-8628740675040476216"
    80 = "2879104050616670826"
    81 = "7874896233432519209"
    82 = "-8914059309599259301"
    83 = "-6848150808824375545"
    84 = "-8136266109847932246"
    85 = "-802813319062420830"
    86 = "-6736413499540721321"
    87 = "3917367986763808549"
    88 = "-379666681410107570"
    89 = "7882932941196884806"
    90 = "-7882385328087490693"
    91 = "4026045856970347435"
    92 = "6023417813053480307"
    93 = "1454228011034335609"
    94 = "5025001820813646212"
    95 = "6499918952172172178"
    96 = "9012395923092219144"
}

foreach ($row in $docIds.Keys) {
    $cell = $ws.Range("L$row")
    $cell.NumberFormat = "@"
    $cell.Value = $docIds[$row]
    $cell.ClearFormats()
}
